# Update workbook with quick lock codes.

$wb = $excel.ActiveWorkbook

# --- Worksheets ---
$wsConnectors = $wb.Worksheets.Item("Connectors")
$wsBehaviour  = $wb.Worksheets.Item("Behaviour Area")

# --- Insert two new rows into the Connectors sheet, above the existing ---
# --- "Profile end cap" row (currently row 6), pushing rows 6-8 to 8-10. ---
$wsConnectors.Rows.Item(6).Insert()
$wsConnectors.Rows.Item(6).Insert()

# --- New row 5: Quick Lock (E slot) ---
$wsConnectors.Range("A5").Value = "Quick Lock (E slot)"
$wsConnectors.Range("C5").Value = 8
$wsConnectors.Range("D5").Value = "MayTec"
$wsConnectors.Range("E5").Value = "1.32.4F2M4.25"

# --- New row 6: Quick Lock (Bolt) ---
$wsConnectors.Range("A6").Value = "Quick Lock (Bolt)"
$wsConnectors.Range("C6").Value = 8
$wsConnectors.Range("D6").Value = "MayTec"
$wsConnectors.Range("E6").Value = "1.64.5419"
$wsConnectors.Range("F6").Value = "For different panel thicknesses use different Hmax value. Hmax should be > panel thickness"

# --- Update sheet selections / active tab ---
[void]$wsBehaviour.Range("E23").Select()
[void]$wsConnectors.Range("A9").Select()
[void]$wsConnectors.Activate()

$wb.Save()
